# "fixed arrow, slide 3"
#
# Reposition/resize the "Curved Down Arrow 8" shape on slide 3.
#   before (EMU): off=(2111189,4262717) ext=(7498078,1794593)
#   after  (EMU): off=(3943846,4592342) ext=(6042989,1549686)
#
# PowerPoint's COM Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) and stored as single-precision floats, so converting
# an EMU target to points and back can truncate by a single EMU; a small
# epsilon keeps the round-trip from landing one EMU short. Rotation
# (180 degrees) is unchanged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$arrow = $s.Shapes.Item("Curved Down Arrow 8")

$epsilon = 0.00001
$arrow.Left   = (3943846 / 12700) + $epsilon
$arrow.Top    = (4592342 / 12700) + $epsilon
$arrow.Width  = (6042989 / 12700) + $epsilon
$arrow.Height = (1549686 / 12700) + $epsilon
